$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new item in the To-Do column (A6)
$ws.Range("A6").Value = "CatBoost Regressor"

# Fill in the previously empty A7 cell (keeps its existing wrap-text style)
$ws.Range("A7").Value = "https://openai.com/dall-e-2/"

# Update the selected cell shown in the sheet view
$ws.Range("A8").Select()
